$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove row 17 (the "Verify field validation in create account" test row)
$ws.Rows.Item(17).Delete()

# Update Execute column (B) from "No" to "Yes" for rows 2-9
$ws.Range("B2:B9").Value = "Yes"

# Update row 16 (now "Verify Create Account with Invalid Data") - mark as Execute = Yes,
# point to the renamed test-data sheet, and use the corrected parameter name
$ws.Range("B16").Value = "Yes"
$ws.Range("C16").Value = "testdata.xls,SignupPage"
$ws.Range("H16").Value = "coyni_mobile.tests.SignUpTest,`ntestCreateAccountInvalidData,`n-pcreateAccount,`n-pfirstName,`n-plastName,`n-pemail,`n-pphoneNumber,`n-ppassword,`n-pconfirmPassword,`n-perrMessage,`n-pelementName"
